$d = $word.ActiveDocument

function Replace-UniqueText($OldText, $NewText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
    if ($found) {
        $rng.Text = $NewText
    } else {
        Write-Output "NOT FOUND: $OldText"
    }
}

# --- Ambiguous heading: 'Political Research and Data Analysis' appears twice.
# Only the occurrence under 'PROGRAMMER - Lake Research Partners' (the second
# one in document order) changes; the one under 'RESEARCH DIRECTOR' stays as-is.
$targetText = 'Political Research and Data Analysis'
$newHeading = 'Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns'
$matchCount = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText) {
        $matchCount = $matchCount + 1
        if ($matchCount -eq 2) {
            $r = $p.Range
            $r.End = $r.End - 1
            $r.Text = $newHeading
        }
    }
}

# --- Remaining unambiguous replacements
Replace-UniqueText '• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations' '• Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions'
Replace-UniqueText '• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics' '• Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics and demographic analysis'
Replace-UniqueText '• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets' '• Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets'
Replace-UniqueText '• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering' '• Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering'
Replace-UniqueText '• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications' '• Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications'
Replace-UniqueText '• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices' '• Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices'
Replace-UniqueText '• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES' '• Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES'
Replace-UniqueText '• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions' '• Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions'
Replace-UniqueText '• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI' '• Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI'
Replace-UniqueText '• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company''s distinguishing products' '• Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company''s distinguishing products'
Replace-UniqueText '• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices' '• Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices'
Replace-UniqueText '• Developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies' '• Architected and developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies'
Replace-UniqueText '• Liaised with officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to improve platform using Docker and Ubuntu' '• Collaborated with senior officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to enhance platform using Docker and Ubuntu'
Replace-UniqueText '• Conceived and built application using Python, Pandas, and Jupyter to predict how crisis economies respond to different humanitarian interventions' '• Conceived and developed predictive application using Python, Pandas, and Jupyter to forecast how crisis economies respond to different humanitarian interventions'
Replace-UniqueText '• Developed RACSO, a web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting' '• Architected and developed RACSO, a comprehensive web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting'
Replace-UniqueText '• Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner' '• Led RFP process and analyzed bids from 1,200 vendors before selecting optimal implementation partner'
Replace-UniqueText '• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research' '• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions'
Replace-UniqueText '• Designed survey deployment system facilitating thousands of simultaneous phone surveys' '• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly $1 million annually in polling costs'
Replace-UniqueText '• Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system' '• Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system'
Replace-UniqueText '• Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill' '• Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill'
Replace-UniqueText '• Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization' '• Led technology operations for multi-million dollar organization while assisting in search for full-time CTO'
Replace-UniqueText '• Made all technology decisions and practices for massive multinational non-governmental organization' '• Directed all technology decisions and practices for massive multinational non-governmental organization'
Replace-UniqueText '• Wrote comprehensive frameworks for internal and external technology audits' '• Developed comprehensive frameworks for internal and external technology audits'
Replace-UniqueText '• Trained beneficiaries on spatial and Census data analysis for public health research' '• Led training initiatives for beneficiaries on spatial and Census data analysis for public health research'
Replace-UniqueText '• Trained NGO staff in web development using Drupal, PHP, and MySQL' '• Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL'
Replace-UniqueText '• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections' '• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions'
Replace-UniqueText 'Political Field Operations and Data Management' 'Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns'
Replace-UniqueText '• Administered all quantitative and qualitative research operations ensuring reporting accuracy' '• Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions'
Replace-UniqueText '• Managed comprehensive survey fielding for multi-million dollar research firm' '• Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm'
